{"js": "// The three <id> values in the transcription need their \"aN\" suffix\n// normalized to plain \"N\" (e.g. \"p089r_a1\" -> \"p089r_1\"), matching the\n// newly downloaded tc/tcn/tl ids. Each <id>...</id> snippet is split\n// across three runs (the literal \"<id>\" tag, the id text itself, and the\n// literal \"</id>\" tag) with differing formatting; searching for/replacing\n// the whole \"<id>oldId</id>\" span collapses it back into a single run\n// (taking on the formatting of the first/opening-tag run), exactly as the\n// target diff shows.\nconst body = context.document.body;\n\nconst idRenames = [\n  [\"p089r_a1\", \"p089r_1\"],\n  [\"p089r_a2\", \"p089r_2\"],\n  [\"p089r_a3\", \"p089r_3\"],\n];\n\nfor (const [oldId, newId] of idRenames) {\n  const searchResults = body.search(`<id>${oldId}</id>`, { matchCase: true });\n  searchResults.load(\"items\");\n  await context.sync();\n\n  for (const range of searchResults.items) {\n    range.insertText(`<id>${newId}</id>`, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# The three <id> values in the transcription need their \"aN\" suffix\n# normalized to plain \"N\" (e.g. \"p089r_a1\" -> \"p089r_1\"), matching the\n# newly downloaded tc/tcn/tl ids. Each <id>...</id> snippet is split\n# across three runs (the literal \"<id>\" tag, the id text itself, and the\n# literal \"</id>\" tag) with differing formatting; finding/replacing the\n# whole \"<id>oldId</id>\" span collapses it back into a single run (taking\n# on the formatting of the first/opening-tag run), exactly as the target\n# diff shows.\n$d = $word.ActiveDocument\n\n$idRenames = @(\n    , @(\"p089r_a1\", \"p089r_1\")\n    , @(\"p089r_a2\", \"p089r_2\")\n    , @(\"p089r_a3\", \"p089r_3\")\n)\n\nforeach ($pair in $idRenames) {\n    $oldId = $pair[0]\n    $newId = $pair[1]\n\n    $range = $d.Content\n    $range.Find.Execute(\n        \"<id>$oldId</id>\",   # FindText\n        $false,              # MatchCase\n        $false,              # MatchWholeWord\n        $false,              # MatchWildcards\n        $false,              # MatchSoundsLike\n        $false,              # MatchAllWordForms\n        $true,               # Forward\n        1,                   # Wrap (wdFindContinue)\n        $false,              # Format\n        \"<id>$newId</id>\",   # ReplaceWith\n        2                    # Replace (wdReplaceOne)\n    ) | Out-Null\n}\n"}
